# Apply the edit described by the diff:
#  - Set the title of the existing slide 2 to "Datasets"
#  - Insert 13 new "Title and Content" slides (positions 3..15) with titles

$p = $ppt.ActivePresentation

# --- Slide 2 ("Datasets") -------------------------------------------------
$s2 = $p.Slides.Item(2)
$s2.Shapes.Item(1).TextFrame.TextRange.Text = "Datasets"

# --- Titles for the new slides, in final presentation order --------------
$titles = @(
    "Preprocessing",
    "Decision Tree Regressor",
    "Random Forest Regressor",
    "Hyperparameter optimization",
    "Comparison to scikit learn DecisionTreeRegressor ",
    "Comparison to scikit learn RandomForestRegressor",
    "Comparison to LLM version",
    "Conclusion",
    "",
    "",
    "",
    "",
    ""
)

$layoutTitleAndContent = 2
$insertPos = 3

for ($i = 0; $i -lt $titles.Count; $i++) {
    $newSlide = $p.Slides.Add($insertPos, $layoutTitleAndContent)
    if ($titles[$i] -ne "") {
        $newSlide.Shapes.Item(1).TextFrame.TextRange.Text = $titles[$i]
    }
    $insertPos = $insertPos + 1
}
